$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "m312"
$ws.Range("E3").Value = "m312"
$ws.Range("E4").Value = "m312"
